$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Medium Ambiguity"
$ws.Range("A4").Value = "Medium Ambiguity"
$ws.Range("A5").Value = "Medium Ambiguity"
$ws.Range("A6").Value = "Medium Ambiguity"
$ws.Range("A8").Value = "Low Ambiguity"
